# water sampling titrations 9/27/2021
# Fill in the two new CRM-accuracy rows (119 existing placeholder row, 120 new row)
# at the bottom of the "Sheet1" data table, then extend the "% off" shared
# formula down through both rows and move the selection to A121 (the next
# empty row) to match Excel's post-entry cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 119 (already existed as an empty placeholder row with A119 formatted
#     as a date) -- fill in the sample data ---------------------------------
$ws.Range("A119").Value = 44466
$ws.Range("B119").Value = 2265.0046219237502
$ws.Range("C119").Value = 2230.52
$ws.Range("E119").Value = 183
$ws.Range("F119").Value = "CRM opened 9/24/2021"

# --- Row 120 (brand-new row) -- copy the date format from row 118 first so
#     A120 keeps the same date number format/style as the rows above it,
#     then fill in the values --------------------------------------------
$ws.Range("A118").Copy()
$ws.Range("A120").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A120").Value = 44466
$ws.Range("B120").Value = 2245.7404688474599
$ws.Range("C120").Value = 2230.52
$ws.Range("E120").Value = 183
$ws.Range("F120").Value = "CRM opened 9/24/2021"

# --- "% off" formula (column D) for both new rows, continuing the same
#     shared formula pattern used by D116:D118 -----------------------------
$ws.Range("D119:D120").Formula = "=100*(B119-C119)/C119"

# --- Move the active selection to the next empty row, as Excel would after
#     data entry -------------------------------------------------------------
$ws.Range("A121").Select()
